# Applies a reordering of tied-rank categories (same underlying value,
# so the raw numbers in column B stay the same) on two sheets:
#   "max-arrecad" : swap A6<->A7, A9<->A10, A17<->A18
#   "tx-sucesso"  : swap A15<->A16

$wb = $excel.ActiveWorkbook

$wsMax = $wb.Worksheets.Item("max-arrecad")
$wsTx  = $wb.Worksheets.Item("tx-sucesso")

function Swap-CellValues($ws, $cellA, $cellB) {
    $valA = $ws.Range($cellA).Value2
    $valB = $ws.Range($cellB).Value2
    $ws.Range($cellA).Value = $valB
    $ws.Range($cellB).Value = $valA
}

# max-arrecad sheet swaps
Swap-CellValues $wsMax "A6" "A7"
Swap-CellValues $wsMax "A9" "A10"
Swap-CellValues $wsMax "A17" "A18"

# tx-sucesso sheet swap
Swap-CellValues $wsTx "A15" "A16"
